# Version 9.3 Fixed Player ratings

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# --- Row 4: 90min Adjuster ---
$ws.Range("B4").Value = -540
$ws.Range("C4").Value = -700

# --- Row 11: Goals minus xG times 100 (only if negative) ---
$row11Cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
foreach ($col in $row11Cols) {
    $ws.Range($col + "11").Value = -3
}

# --- Row 19: Passes Incomplete ---
$row19 = @{
    "B19" = 50;  "C19" = -19; "D19" = -19; "E19" = -19; "F19" = 50;
    "G19" = 50;  "H19" = 50;  "I19" = 50;  "J19" = 50;  "K19" = 50;
    "L19" = 50;  "M19" = 50;  "N19" = -19; "O19" = 50;  "P19" = 50;
    "Q19" = 50;  "R19" = 50;  "S19" = 50;  "T19" = -19; "U19" = 50;
    "V19" = 50;  "W19" = -19; "X19" = 50;  "Y19" = 50;  "Z19" = 50;
    "AA19" = 50; "AB19" = -19
}
foreach ($key in $row19.Keys) {
    $ws.Range($key).Value = $row19[$key]
}

# --- Row 26: Open Play xA ---
$row26 = @{
    "B26" = 50;  "C26" = 50;  "D26" = 50;  "E26" = 10;  "F26" = 50;
    "G26" = 50;  "H26" = 10;  "I26" = 10;  "J26" = 10;  "K26" = 10;
    "L26" = 10;  "M26" = 50;  "N26" = 50;  "O26" = 10;  "P26" = 10;
    "Q26" = 10;  "R26" = 10;  "S26" = 10;  "T26" = 50;  "U26" = 10;
    "V26" = 10;  "W26" = 50;  "X26" = 10;  "Y26" = 10;  "Z26" = 50;
    "AA26" = 10; "AB26" = 50
}
foreach ($key in $row26.Keys) {
    $ws.Range($key).Value = $row26[$key]
}

# --- Row 39: Saves ---
$ws.Range("B39").Value = 120

# --- Column A width (best-fit width for the longest role name) ---
$ws.Columns.Item(1).ColumnWidth = 39

# --- Sheet view: scroll + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B40").Select()
